# The deck's theme (ppt/theme/theme1.xml, bound to the one Slide Master)
# currently carries the "Integral" colour scheme. The target edit swaps it
# for the stock "Office Theme" colour scheme (the 12 dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink colours), matching the theme that was otherwise only sitting
# unused in ppt/theme/theme2.xml (referenced solely by the Notes Master).
#
# PowerPoint's COM model exposes these 12 theme colours via
# Master.Theme.ThemeColorScheme(1..12).RGB — setting them rewrites the
# <a:clrScheme> of the theme part backing the Slide Master.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$tcs = $m.Theme.ThemeColorScheme

# Office Theme colour scheme, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$tcs.Item(1).RGB  = 0x000000   # dk1      000000
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = 0x6A5444   # dk2      44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink 954F72
